$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column C ("Förändrad") on every data row: 45184 -> 45186 ---------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45186
    }
}

# --- 2) Hyperlink formulas in columns S, T, V, W, X, Y: add the display
#        text (the "Beteckning" in column A) as the second HYPERLINK arg --
$linkCols = @("S", "T", "V", "W", "X", "Y")

for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Range("A" + $r).Value2
    if ([string]::IsNullOrEmpty($label)) { continue }

    foreach ($col in $linkCols) {
        $rng = $ws.Range($col + $r)
        $f = $rng.Formula
        if ([string]::IsNullOrEmpty($f)) { continue }
        if ($f -notlike '*HYPERLINK(*') { continue }
        if ($f -like '*,*') { continue }  # already has a second argument

        $trimmed = $f.TrimEnd()
        if ($trimmed.EndsWith(")")) {
            $newFormula = $trimmed.Substring(0, $trimmed.Length - 1) + ', "' + $label + '")'
            $rng.Formula = $newFormula
        }
    }
}
